# Commit: "Added new files for equity calculation"
# - Column I (trade_date) for the already-existing rows (2..356) switches
#   from the date-only number format to the datetime number format (the
#   same format already used by column B / datetime).
# - Four new rows (357..360) are appended with the next four trading days,
#   continuing the OLA Electric Mobility Ltd. history. On these new rows,
#   column I keeps using the (older) date-only format that column I used
#   to have before this edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Re-format the existing trade_date column (I2:I356) to the datetime format.
$ws.Range("I2:I356").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# 2) Append the four new trading-day rows (357-360).
$newRows = @(
    @{ Row=357; Close=31.54; Date=46049; High=32.22; Low=30.76; Open=32.15; Volume=65726680 },
    @{ Row=358; Close=32.08; Date=46050; High=32.47; Low=31.66; Open=32.07; Volume=54665824 },
    @{ Row=359; Close=31.85; Date=46051; High=32.13; Low=31.31; Open=32.13; Volume=39698821 },
    @{ Row=360; Close=32.33; Date=46052; High=32.78; Low=31.42; Open=31.62; Volume=63348957 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.Close
    $ws.Cells.Item($row, 2).Value = $r.Date
    $ws.Cells.Item($row, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 3).Value = "NSE"
    $ws.Cells.Item($row, 4).Value = $r.High
    $ws.Cells.Item($row, 5).Value = $r.Low
    $ws.Cells.Item($row, 6).Value = $r.Open
    $ws.Cells.Item($row, 7).Value = "OLAELE"
    $ws.Cells.Item($row, 8).Value = $r.Volume
    $ws.Cells.Item($row, 9).Value = $r.Date
    $ws.Cells.Item($row, 9).NumberFormat = "YYYY-MM-DD"
    $ws.Cells.Item($row, 10).Value = "INE0LXG01040"
    $ws.Cells.Item($row, 11).Value = "OLA Electric Mobility Ltd"
    $ws.Cells.Item($row, 12).Value = "OLAELE"
    $ws.Cells.Item($row, 13).Value = "BREEZE"
}
